$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (H = "home"/team row 2 values), divisional round update
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 310
$wsOff.Range("C2").Value = 220
$wsOff.Range("D2").Value = 77
$wsOff.Range("E2").Value = 38
$wsOff.Range("G2").Value = 2

# DEF sheet - row 2
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 281
$wsDef.Range("C2").Value = 198
$wsDef.Range("D2").Value = 67
$wsDef.Range("E2").Value = 32
